$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change L1's value from "WS" to a new shared string "WA"
$ws.Range("L1").Value2 = "WA"

# Update the selected range / active cell shown in the sheet view
$ws.Range("S12").Select()
